# Auto-generated edit script: rebuild rows 2-36 of the worksheet per the target diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "Rv3423c"
$ws.Cells.Item(2, 2).Value = 5
$ws.Cells.Item(2, 3).Value = "alr Rv3423c MTCY78.06"
$ws.Cells.Item(2, 4).Value = "FUNCTION: Catalyzes the interconversion of L-alanine and D-alanine. D-alanine plays a key role in peptidoglycan cross-linking. {ECO:0000269|PubMed:11267762}."
$ws.Cells.Item(2, 5).Value = 7

$ws.Cells.Item(3, 1).Value = "Rv1485"
$ws.Cells.Item(3, 2).Value = 5
$ws.Cells.Item(3, 3).Value = "cpfC hemH hemZ Rv1485 MTCY277.06"
$ws.Cells.Item(3, 4).Value = "FUNCTION: Involved in coproporphyrin-dependent heme b biosynthesis (PubMed:25646457). Catalyzes the insertion of ferrous iron into coproporphyrin III to form Fe-coproporphyrin III (PubMed:25646457). Has weaker activity with coproporphyrin I, protoporphyrin IX, deuteroporphyrin, 2,4 hydroxyethyl and 2,4 disulfonate (PubMed:25646457, PubMed:11948160). {ECO:0000269|PubMed:11948160, ECO:0000269|PubMed:25646457}."
$ws.Cells.Item(3, 5).Value = 7

$ws.Cells.Item(4, 1).Value = "Rv1099c"
$ws.Cells.Item(4, 2).Value = 5
$ws.Cells.Item(4, 3).Value = "glpX Rv1099c"
$ws.Cells.Item(4, 4).Value = "FUNCTION: Catalyzes the hydrolysis of fructose 1,6-bisphosphate to fructose 6-phosphate (PubMed:15470127, PubMed:21451980). Seems to be the major FBPase of M.tuberculosis and to play a key role in gluconeogenesis for conversion of lipid carbon into cell wall glycans. Does not display activity against inositol 1-phosphate (PubMed:15470127). {ECO:0000269|PubMed:15470127, ECO:0000269|PubMed:21451980}."
$ws.Cells.Item(4, 5).Value = 7

$ws.Cells.Item(5, 1).Value = "Rv2097c"
$ws.Cells.Item(5, 2).Value = 5
$ws.Cells.Item(5, 3).Value = "pafA paf Rv2097c MTCY49.37c"
$ws.Cells.Item(5, 4).Value = "FUNCTION: Catalyzes the covalent attachment of the prokaryotic ubiquitin-like protein modifier Pup to the proteasomal substrate proteins, thereby targeting them for proteasomal degradation. This tagging system is termed pupylation. The ligation reaction involves the side-chain carboxylate of the C-terminal glutamate of Pup and the side-chain amino group of a substrate lysine. PafA is required to confer resistance against the lethal effects of reactive nitrogen intermediates (RNI), antimicrobial molecules produced by activated macrophages and other cell types. {ECO:0000269|PubMed:14671303, ECO:0000269|PubMed:17082771, ECO:0000269|PubMed:19448618, ECO:0000269|PubMed:20355727}."
$ws.Cells.Item(5, 5).Value = 7

$ws.Cells.Item(6, 1).Value = "Rv1854c"
$ws.Cells.Item(6, 2).Value = 5
$ws.Cells.Item(6, 3).Value = "ndh Rv1854c"
$ws.Cells.Item(6, 4).Value = "FUNCTION: Alternative, nonproton pumping NADH:quinone oxidoreductase that delivers electrons to the respiratory chain by oxidation of NADH and reduction of quinones (PubMed:15767566, PubMed:29382761, PubMed:29522317). Ndh is probably the main NADH dehydrogenase of M.tuberculosis (PubMed:29382761). {ECO:0000269|PubMed:15767566, ECO:0000269|PubMed:29382761, ECO:0000269|PubMed:29522317}."
$ws.Cells.Item(6, 5).Value = 7

$ws.Cells.Item(7, 1).Value = "Rv0998"
$ws.Cells.Item(7, 2).Value = 5
$ws.Cells.Item(7, 3).Value = "Rv0998"
$ws.Cells.Item(7, 4).Value = "FUNCTION: Catalyzes specifically the acetylation of the epsilon-amino group of a highly conserved lysine residue in acetyl-CoA synthetase (ACS). This acetylation results in the inactivation of ACS activity and could be important for mycobacteria to adjust to environmental changes. {ECO:0000269|PubMed:20507997, ECO:0000269|PubMed:21627103, ECO:0000269|PubMed:22773105}."
$ws.Cells.Item(7, 5).Value = 7

$ws.Cells.Item(8, 1).Value = "Rv3543c"
$ws.Cells.Item(8, 2).Value = 5
$ws.Cells.Item(8, 3).Value = "fadE29 Rv3543c"
$ws.Cells.Item(8, 4).Value = "FUNCTION: Involved in the third cycle of side chain dehydrogenation in the beta-oxidation of cholesterol catabolism (PubMed:26161441). Contributes partly to the virulence by increasing the efficiency of beta-oxidation (PubMed:22045806, PubMed:23560677). Catalyzes the dehydrogenation of 2'-propanoyl-CoA ester side chains of 3-oxo-4-pregnene-20-carboxyl-CoA (3-OPC-CoA) to yield 3-oxo-4,17-pregnadiene-20-carboxyl-CoA (3-OPDC-CoA). Also able to dehydrogenate steroyl-CoA such as 3-oxo-chol-4-en-24-oyl-CoA (3-OCO-CoA), 1beta-(2'-propanoyl-CoA)-3a-alpha-H- 7a-beta-methylhexahydro-4-indanone (indanone-CoA ester), hexahydroindanone and pregenenone (PubMed:22045806, PubMed:23560677). {ECO:0000269|PubMed:22045806, ECO:0000269|PubMed:23560677, ECO:0000269|PubMed:26161441}."
$ws.Cells.Item(8, 5).Value = 7

$ws.Cells.Item(9, 1).Value = "Rv0956"
$ws.Cells.Item(9, 2).Value = 5
$ws.Cells.Item(9, 3).Value = "purN Rv0956"
$ws.Cells.Item(9, 4).Value = "FUNCTION: Catalyzes the transfer of a formyl group from 10-formyltetrahydrofolate to 5-phospho-ribosyl-glycinamide (GAR), producing 5-phospho-ribosyl-N-formylglycinamide (FGAR) and tetrahydrofolate. {ECO:0000255|HAMAP-Rule:MF_01930}."
$ws.Cells.Item(9, 5).Value = 7

$ws.Cells.Item(10, 1).Value = "Rv2115c"
$ws.Cells.Item(10, 2).Value = 5
$ws.Cells.Item(10, 3).Value = "mpa Rv2115c MTCY261.11c"
$ws.Cells.Item(10, 4).Value = "FUNCTION: ATPase which is responsible for recognizing, binding, unfolding and translocation of pupylated proteins into the bacterial 20S proteasome core particle. May be essential for opening the gate of the 20S proteasome via an interaction with its C-terminus, thereby allowing substrate entry and access to the site of proteolysis. Thus, the C-termini of the proteasomal ATPase may function like a 'key in a lock' to induce gate opening and therefore regulate proteolysis. Is required but not sufficient to confer resistance against the lethal effects of reactive nitrogen intermediates (RNI), antimicrobial molecules produced by activated macrophages and other cell types. {ECO:0000255|HAMAP-Rule:MF_02112, ECO:0000269|PubMed:14671303, ECO:0000269|PubMed:15659170, ECO:0000269|PubMed:17082771, ECO:0000269|PubMed:19836337, ECO:0000269|PubMed:20203624}."
$ws.Cells.Item(10, 5).Value = 7

$ws.Cells.Item(11, 1).Value = "Rv0350"
$ws.Cells.Item(11, 2).Value = 5
$ws.Cells.Item(11, 3).Value = "dnaK Rv0350 MTCY13E10.10"
$ws.Cells.Item(11, 4).Value = "FUNCTION: Acts as a chaperone. {ECO:0000255|HAMAP-Rule:MF_00332}.; FUNCTION: Recombinant extracellular protein activates expression of NF-kappa-B in immortalized human dermal endothelial cells in a TLR2- and TLR4-dependent manner. Activation occurs via MYD88-dependent and -independent pathways and requires TIRAP, TRIF and TRAM (some experiments done in mouse cells, mice do not usually catch tuberculosis) (PubMed:15809303). {ECO:0000269|PubMed:15809303}."
$ws.Cells.Item(11, 5).Value = 7

$ws.Cells.Item(12, 1).Value = "Rv1712"
$ws.Cells.Item(12, 2).Value = 5
$ws.Cells.Item(12, 3).Value = "cmk Rv1712 MTCI125.34"
$ws.Cells.Item(12, 4).Value = ""
$ws.Cells.Item(12, 5).Value = 7

$ws.Cells.Item(13, 1).Value = "Rv0440"
$ws.Cells.Item(13, 2).Value = 5
$ws.Cells.Item(13, 3).Value = "groEL2 groEL-2 groL2 hsp65 mtc28 Rv0440 MTV037.04"
$ws.Cells.Item(13, 4).Value = "FUNCTION: Prevents aggregation of substrate proteins and promotes their refolding (PubMed:15327959). {ECO:0000269|PubMed:15327959}.; FUNCTION: Recombinant extracellular protein activates expression of NF-kappa-B in immortalized human dermal endothelial cells in a TLR4-dependent, TLR2-independent manner. Activation occurs via MYD88-dependent and -independent pathways and requires TIRAP, TRIF, TRAM and MD-2 (some experiments done in mouse cells, mice do not usually catch tuberculosis) (PubMed:15809303). {ECO:0000269|PubMed:15809303}."
$ws.Cells.Item(13, 5).Value = 7

$ws.Cells.Item(14, 1).Value = "Rv2881c"
$ws.Cells.Item(14, 2).Value = 4
$ws.Cells.Item(14, 3).Value = "cdsA Rv2881c MTCY274.12c"
$ws.Cells.Item(14, 4).Value = ""
$ws.Cells.Item(14, 5).Value = 7

$ws.Cells.Item(15, 1).Value = "Rv1655"
$ws.Cells.Item(15, 2).Value = 4
$ws.Cells.Item(15, 3).Value = "argD Rv1655 MTCY06H11.20"
$ws.Cells.Item(15, 4).Value = ""
$ws.Cells.Item(15, 5).Value = 7

$ws.Cells.Item(16, 1).Value = "Rv2461c"
$ws.Cells.Item(16, 2).Value = 4
$ws.Cells.Item(16, 3).Value = "clpP1 clpP Rv2461c MTV008.17c"
$ws.Cells.Item(16, 4).Value = "FUNCTION: Cleaves peptides in various proteins in a process that requires ATP hydrolysis. Has a chymotrypsin-like activity. Plays a major role in the degradation of misfolded proteins (By similarity). Degrades anti-sigma-D factor (rsdA) when present in a complex with ClpP2 and ClpX. Does not seem to act on anti-sigma-L factor (rslA). {ECO:0000255|HAMAP-Rule:MF_00444, ECO:0000269|PubMed:23314154}."
$ws.Cells.Item(16, 5).Value = 7

$ws.Cells.Item(17, 1).Value = "Rv2196"
$ws.Cells.Item(17, 2).Value = 4
$ws.Cells.Item(17, 3).Value = "qcrB Rv2196 MTCY190.07"
$ws.Cells.Item(17, 4).Value = "FUNCTION: Cytochrome b subunit of the cytochrome bc1 complex, an essential component of the respiratory electron transport chain required for ATP synthesis. The bc1 complex catalyzes the oxidation of ubiquinol and the reduction of cytochrome c in the respiratory chain. The bc1 complex operates through a Q-cycle mechanism that couples electron transfer to generation of the proton gradient that drives ATP synthesis. The cytochrome b subunit contains two ubiquinol reactive sites: the oxidation (QP) site and the reduction (QN) site. {ECO:0000305, ECO:0000305|PubMed:23913123, ECO:0000305|PubMed:26158909}."
$ws.Cells.Item(17, 5).Value = 7

$ws.Cells.Item(18, 1).Value = "Rv1110"
$ws.Cells.Item(18, 2).Value = 4
$ws.Cells.Item(18, 3).Value = "ispH2 lytB2 Rv1110 MTV017.63"
$ws.Cells.Item(18, 4).Value = "FUNCTION: Catalyzes the conversion of 1-hydroxy-2-methyl-2-(E)-butenyl 4-diphosphate (HMBPP) into a mixture of isopentenyl diphosphate (IPP) and dimethylallyl diphosphate (DMAPP) (PubMed:23091471). Acts in the terminal step of the DOXP/MEP pathway for isoprenoid precursor biosynthesis. Has a higher activity compared with LytB2 (PubMed:23091471). Is essential for M.tuberculosis growth in vitro (PubMed:26309039). {ECO:0000255|HAMAP-Rule:MF_00191, ECO:0000269|PubMed:23091471, ECO:0000269|PubMed:26309039}."
$ws.Cells.Item(18, 5).Value = 7

$ws.Cells.Item(19, 1).Value = "Rv1796"
$ws.Cells.Item(19, 2).Value = 3
$ws.Cells.Item(19, 3).Value = "mycP5 Rv1796 LH57_09820"
$ws.Cells.Item(19, 4).Value = ""
$ws.Cells.Item(19, 5).Value = 7

$ws.Cells.Item(20, 1).Value = "Rv0286"
$ws.Cells.Item(20, 2).Value = 3
$ws.Cells.Item(20, 3).Value = "PPE4 Rv0286"
$ws.Cells.Item(20, 4).Value = "FUNCTION: Important for the siderophore-mediated iron-acquisition function of ESX-3. {ECO:0000269|PubMed:26729876}."
$ws.Cells.Item(20, 5).Value = 7

$ws.Cells.Item(21, 1).Value = "Rv3540c"
$ws.Cells.Item(21, 2).Value = 3
$ws.Cells.Item(21, 3).Value = "ltp2 Rv3540c"
$ws.Cells.Item(21, 4).Value = "FUNCTION: Involved in cholesterol side chain degradation (PubMed:22045806, PubMed:29109182). When associated with the ChsH1/ChsH2 hydratase, catalyzes the retroaldol cleavage of 17-hydroxy-3-oxo-4-pregnene-20-carboxyl-CoA (17-HOPC-CoA) produced by the hydratase, forming androst-4-ene-3,17-dione and propionyl-CoA (PubMed:29109182, PubMed:31568719). {ECO:0000269|PubMed:22045806, ECO:0000269|PubMed:29109182, ECO:0000269|PubMed:31568719}."
$ws.Cells.Item(21, 5).Value = 7

$ws.Cells.Item(22, 1).Value = "Rv1340"
$ws.Cells.Item(22, 2).Value = 3
$ws.Cells.Item(22, 3).Value = "rph rphA Rv1340 MTCY02B10.04 MTCY130.25"
$ws.Cells.Item(22, 4).Value = "FUNCTION: Phosphorolytic 3'-5' exoribonuclease that plays an important role in tRNA 3'-end maturation. Removes nucleotide residues following the 3'-CCA terminus of tRNAs; can also add nucleotides to the ends of RNA molecules by using nucleoside diphosphates as substrates, but this may not be physiologically important. Probably plays a role in initiation of 16S rRNA degradation (leading to ribosome degradation) during starvation. {ECO:0000255|HAMAP-Rule:MF_00564}."
$ws.Cells.Item(22, 5).Value = 7

$ws.Cells.Item(23, 1).Value = "Rv0957"
$ws.Cells.Item(23, 2).Value = 3
$ws.Cells.Item(23, 3).Value = "purH Rv0957 MTCY10D7.17c"
$ws.Cells.Item(23, 4).Value = ""
$ws.Cells.Item(23, 5).Value = 7

$ws.Cells.Item(24, 1).Value = "Rv2977c"
$ws.Cells.Item(24, 2).Value = 3
$ws.Cells.Item(24, 3).Value = "thiL Rv2977c"
$ws.Cells.Item(24, 4).Value = "FUNCTION: Catalyzes the ATP-dependent phosphorylation of thiamine-monophosphate (TMP) to form thiamine-pyrophosphate (TPP), the active form of vitamin B1. {ECO:0000255|HAMAP-Rule:MF_02128}."
$ws.Cells.Item(24, 5).Value = 7

$ws.Cells.Item(25, 1).Value = "Rv1166"
$ws.Cells.Item(25, 2).Value = 3
$ws.Cells.Item(25, 3).Value = "lpqW Rv1166"
$ws.Cells.Item(25, 4).Value = "FUNCTION: May directly or indirectly regulate the accessibility of the key branch point intermediate, monoacyl phosphatidylinositol tetramannoside (AcPIM4), to the elongating alpha-1,6 mannosyltransferases which could regulate the lipoarabinomannans (LAMs) biosynthesis. {ECO:0000250}."
$ws.Cells.Item(25, 5).Value = 7

$ws.Cells.Item(26, 1).Value = "Rv2448c"
$ws.Cells.Item(26, 2).Value = 3
$ws.Cells.Item(26, 3).Value = "valS Rv2448c MTV008.04c"
$ws.Cells.Item(26, 4).Value = "FUNCTION: Catalyzes the attachment of valine to tRNA(Val). As ValRS can inadvertently accommodate and process structurally similar amino acids such as threonine, to avoid such errors, it has a 'posttransfer' editing activity that hydrolyzes mischarged Thr-tRNA(Val) in a tRNA-dependent manner. {ECO:0000255|HAMAP-Rule:MF_02004}."
$ws.Cells.Item(26, 5).Value = 7

$ws.Cells.Item(27, 1).Value = "Rv2846c"
$ws.Cells.Item(27, 2).Value = 2
$ws.Cells.Item(27, 3).Value = "efpA Rv2846c"
$ws.Cells.Item(27, 4).Value = ""
$ws.Cells.Item(27, 5).Value = 7

$ws.Cells.Item(28, 1).Value = "Rv2362c"
$ws.Cells.Item(28, 2).Value = 2
$ws.Cells.Item(28, 3).Value = "recO Rv2362c MTCY27.18"
$ws.Cells.Item(28, 4).Value = "FUNCTION: Involved in DNA repair and RecF pathway recombination. {ECO:0000250}."
$ws.Cells.Item(28, 5).Value = 7

$ws.Cells.Item(29, 1).Value = "Rv1830"
$ws.Cells.Item(29, 2).Value = 1
$ws.Cells.Item(29, 3).Value = "Rv1830 MTCY1A11.13c"
$ws.Cells.Item(29, 4).Value = ""
$ws.Cells.Item(29, 5).Value = 7

$ws.Cells.Item(30, 1).Value = "Rv0023"
$ws.Cells.Item(30, 2).Value = 1
$ws.Cells.Item(30, 3).Value = "Rv0023 MTCY10H4.23"
$ws.Cells.Item(30, 4).Value = ""
$ws.Cells.Item(30, 5).Value = 7

$ws.Cells.Item(31, 1).Value = "Rv1277"
$ws.Cells.Item(31, 2).Value = 1
$ws.Cells.Item(31, 3).Value = "Rv1277"
$ws.Cells.Item(31, 4).Value = ""
$ws.Cells.Item(31, 5).Value = 7

$ws.Cells.Item(32, 1).Value = "Rv2516c"
$ws.Cells.Item(32, 2).Value = 1
$ws.Cells.Item(32, 3).Value = "Rv2516c"
$ws.Cells.Item(32, 4).Value = ""
$ws.Cells.Item(32, 5).Value = 7

$ws.Cells.Item(33, 1).Value = "Rv2535c"
$ws.Cells.Item(33, 2).Value = 1
$ws.Cells.Item(33, 3).Value = "pepQ Rv2535c"
$ws.Cells.Item(33, 4).Value = ""
$ws.Cells.Item(33, 5).Value = 7

$ws.Cells.Item(34, 1).Value = "Rv2382c"
$ws.Cells.Item(34, 2).Value = 1
$ws.Cells.Item(34, 3).Value = "mbtC Rv2382c"
$ws.Cells.Item(34, 4).Value = ""
$ws.Cells.Item(34, 5).Value = 7

$ws.Cells.Item(35, 1).Value = "Rv0238"
$ws.Cells.Item(35, 2).Value = 1
$ws.Cells.Item(35, 3).Value = "Rv0238"
$ws.Cells.Item(35, 4).Value = ""
$ws.Cells.Item(35, 5).Value = 7

$ws.Cells.Item(36, 1).Value = "Rv0358"
$ws.Cells.Item(36, 2).Value = 1
$ws.Cells.Item(36, 3).Value = "Rv0358"
$ws.Cells.Item(36, 4).Value = ""
$ws.Cells.Item(36, 5).Value = 7
